# Rename the header cells from the "_old"/"_new" suffix naming scheme to the
# new "_FV2310"/"_FV2404" naming scheme (column K "diff" is left untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"

$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# Stash a copy of the (about to be re-created) header formatting far away from
# the used range so it survives the table creation below untouched, then
# strip the formatting from the header row. A freshly created ListObject
# derives a header dxf from whatever formatting differences exist on the
# header row at creation time - by making the header "plain" first we avoid
# Excel manufacturing a dxf / headerRowDxfId that the source workbook does
# not have.
$hdr = $ws.Range("A1:U1")
$scratch = $ws.Range("A1000:U1000")
$hdr.Copy($scratch)
$hdr.ClearFormats()

# Turn the data range into an actual Excel Table (ListObject), matching the
# newly introduced xl/tables/table1.xml part.
$rng = $ws.Range("A1:U66")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Restore the original header formatting (bold, shaded fill, centered,
# wrapped, bordered) from the stashed copy, then remove the scratch copy
# again so it doesn't affect the sheet's used range / dimension.
$scratch.Copy()
$hdr.PasteSpecial(-4122)
$scratch.Clear()

# Freeze the header row (pane split after row 1) like the updated sheetView.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
